# The document's footers each anchor a small "Britam Public"
# classification text box (one per footer type: primary, first page,
# even page). The edit removes the visible "Britam Public" caption text
# from each of those text boxes, leaving the (now empty) text box shapes
# in place.
$d = $word.ActiveDocument

$footerTypes = 1, 2, 3   # wdHeaderFooterPrimary, wdHeaderFooterFirstPage, wdHeaderFooterEvenPages

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    foreach ($ft in $footerTypes) {
        try {
            $ftr = $sec.Footers($ft)
        } catch {
            continue
        }

        for ($i = 1; $i -le $ftr.Shapes.Count; $i++) {
            $shp = $ftr.Shapes.Item($i)
            $txt = $shp.TextFrame.TextRange.Text

            if ($txt -like "*Britam Public*") {
                $shp.TextFrame.TextRange.Text = ""
            }
        }
    }
}

Write-Output "Cleared Britam Public classification text boxes in footers"
